# Update cryptocurrency price (D) and volume change (E) cells to match the
# refreshed data feed. Values are forced to text (leading apostrophe) so
# Excel does not reinterpret numeric-looking strings (e.g. "212.32",
# "0.0617") as actual numbers, matching the original inline-string cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.351.99"
$ws.Range("E2").Value = "'  +1.32%  "
$ws.Range("D3").Value = "'1.622.13"
$ws.Range("E3").Value = "'  +1.84%  "
$ws.Range("E4").Value = "'  -0.26%  "
$ws.Range("D5").Value = "'212.32"
$ws.Range("E5").Value = "'  +0.71%  "
$ws.Range("E6").Value = "'  -0.16%  "
$ws.Range("D7").Value = "'0.485"
$ws.Range("E7").Value = "'  +1.18%  "
$ws.Range("E8").Value = "'  +0.42%  "
$ws.Range("D9").Value = "'0.0617"
$ws.Range("E9").Value = "'  +0.76%  "
$ws.Range("D10").Value = "'18.88"
$ws.Range("E10").Value = "'  +4.76%  "
$ws.Range("D11").Value = "'0.0815"
$ws.Range("E11").Value = "'  +0.89%  "
$ws.Range("D12").Value = "'1.844.73"
$ws.Range("E12").Value = "'  +1.45%  "
$ws.Range("D13").Value = "'1.635.47"
$ws.Range("E13").Value = "'  +2.86%  "
$ws.Range("E14").Value = "'  +0.05%  "
$ws.Range("D15").Value = "'0.520"
$ws.Range("E15").Value = "'  +1.40%  "
$ws.Range("D16").Value = "'26.342.18"
$ws.Range("E16").Value = "'  +1.04%  "
$ws.Range("D17").Value = "'62.57"
$ws.Range("E17").Value = "'  +1.77%  "
$ws.Range("E18").Value = "'  +0.52%  "
$ws.Range("E19").Value = "'  -0.09%  "
$ws.Range("D20").Value = "'202.29"
$ws.Range("E20").Value = "'  +0.03%  "
$ws.Range("E21").Value = "'  +0.70%  "
$ws.Range("D22").Value = "'9.35"
$ws.Range("E22").Value = "'  +1.00%  "
$ws.Range("E23").Value = "'  +0.86%  "
$ws.Range("D24").Value = "'1.89"
$ws.Range("E24").Value = "'  -2.64%  "
$ws.Range("D25").Value = "'144.53"
$ws.Range("E25").Value = "'  +0.55%  "
$ws.Range("E26").Value = "'  -0.38%  "
$ws.Range("D27").Value = "'0.120"
$ws.Range("E27").Value = "'  -1.10%  "
$ws.Range("D28").Value = "'15.19"
$ws.Range("E28").Value = "'  -0.08%  "
$ws.Range("D30").Value = "'0.0514"
$ws.Range("E30").Value = "'  +8.15%  "
$ws.Range("E31").Value = "'  +0.47%  "
$ws.Range("E32").Value = "'  +2.24%  "
$ws.Range("D33").Value = "'2.93"
$ws.Range("E33").Value = "'  +1.10%  "
$ws.Range("E34").Value = "'  +1.16%  "
$ws.Range("E35").Value = "'  +2.26%  "
$ws.Range("D36").Value = "'1.177.19"
$ws.Range("E36").Value = "'  +4.22%  "
$ws.Range("D37").Value = "'0.0164"
$ws.Range("E37").Value = "'  +1.04%  "
$ws.Range("D38").Value = "'0.810"
$ws.Range("E38").Value = "'  +2.34%  "
$ws.Range("E39").Value = "'  +0.09%  "
$ws.Range("E40").Value = "'  +0.31%  "
$ws.Range("D41").Value = "'0.498"
$ws.Range("E41").Value = "'  +0.80%  "
$ws.Range("D42").Value = "'5.37"
$ws.Range("E42").Value = "'  +4.39%  "
$ws.Range("D43").Value = "'0.786"
$ws.Range("E43").Value = "'  +0.52%  "
$ws.Range("D44").Value = "'1.757.65"
$ws.Range("E44").Value = "'  +1.66%  "
$ws.Range("D45").Value = "'92.65"
$ws.Range("E45").Value = "'  +0.66%  "
$ws.Range("E46").Value = "'  +2.84%  "
$ws.Range("D47").Value = "'53.84"
$ws.Range("E47").Value = "'  +0.28%  "
$ws.Range("E48").Value = "'  +0.66%  "
$ws.Range("E49").Value = "'  +0.67%  "
$ws.Range("E50").Value = "'  -0.38%  "
$ws.Range("D51").Value = "'7.29"
$ws.Range("E51").Value = "'  +1.80%  "
